# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# worksheet with freshly scraped values, mirroring the GitHub Actions
# commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.923.15"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "2.290.20"
$ws.Range("E3").Value = "  +1.89%  "
$ws.Range("E4").Value = "  -0.21%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.90"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +1.12%  "
$ws.Range("D6").Value = "0.638"
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("D7").Value = "75.03"
$ws.Range("E7").Value = "  +6.79%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "0.649"
$ws.Range("E9").Value = "  -2.72%  "
$ws.Range("D10").Value = "39.03"
$ws.Range("E10").Value = "  -0.52%  "
$ws.Range("D11").Value = "0.0984"
$ws.Range("E11").Value = "  +2.39%  "
$ws.Range("D12").Value = "7.48"
$ws.Range("E12").Value = "  -0.99%  "
$ws.Range("D13").Value = "0.106"
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("D14").Value = "2.633.24"
$ws.Range("E14").Value = "  +1.82%  "
$ws.Range("D15").Value = "15.13"
$ws.Range("E15").Value = "  +2.36%  "
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.870"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  -1.69%  "
$ws.Range("D17").Value = "2.295.94"
$ws.Range("E17").Value = "  +2.18%  "
$ws.Range("D18").Value = "42.800.51"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("E19").Value = "  +1.53%  "
$ws.Range("D20").Value = "6.24"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("D21").Value = "72.32"
$ws.Range("E21").Value = "  -0.92%  "
$ws.Range("D22").Value = "235.85"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").Value = "2.17"
$ws.Range("E23").Value = "  +5.67%  "
$ws.Range("D24").Value = "3.88"
$ws.Range("E24").Value = "  -1.48%  "
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("D26").Value = "11.35"
$ws.Range("E26").Value = "  -1.55%  "
$ws.Range("D27").Value = "2.41"
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("D28").Value = "2.13"
$ws.Range("E28").Value = "  -3.67%  "
$ws.Range("D29").Value = "167.35"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Value = "21.04"
$ws.Range("E30").Value = "  +0.43%  "
$ws.Range("D31").Value = "0.0863"
$ws.Range("E31").Value = "  +9.17%  "
$ws.Range("D32").Value = "6.25"
$ws.Range("E32").Value = "  -3.39%  "
$ws.Range("E33").Value = "  -0.83%  "
$ws.Range("D34").Value = "31.39"
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("D35").Value = "0.127"
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("E36").Value = "  +6.60%  "
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.80"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  +1.97%  "
$ws.Range("D38").Value = "0.0305"
$ws.Range("E38").Value = "  -4.81%  "
$ws.Range("D39").Value = "13.63"
$ws.Range("E39").Value = "  +8.52%  "
$ws.Range("E40").Value = "  -0.35%  "
$ws.Range("D41").Value = "5.97"
$ws.Range("E42").Value = "  +4.86%  "
$ws.Range("E43").Value = "  +1.97%  "
$ws.Range("D44").Value = "61.18"
$ws.Range("E44").Value = "  -2.46%  "
$ws.Range("D45").Value = "4.84"
$ws.Range("E45").Value = "  -0.39%  "
$ws.Range("D46").Value = "105.29"
$ws.Range("E46").Value = "  +11.33%  "
$ws.Range("E47").Value = "  -1.78%  "
$ws.Range("E48").Value = "  +0.18%  "
$ws.Range("D49").Value = "1.16"
$ws.Range("E49").Value = "  -0.76%  "
$ws.Range("E50").Value = "  -1.25%  "
$ws.Range("E51").Value = "  -1.48%  "
